$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text updates -------------------------------------------------
# "Volume 30   Number  14" -> "Volume 30   Number  15"
$ws.Range("A8").Characters(21, 2).Text = "15"

# "Report Covering the Week  4/3/2023  Through  4/9/2023"
#   -> "Report Covering the Week  4/10/2023  Through  4/16/2023"
$ws.Range("C9").Characters(27, 8).Text = "4/10/2023"
$ws.Range("C9").Characters(47, 8).Text = "4/16/2023"

# --- Crime-statistics table updates --------------------------------------
# Donor cells used to carry over the correct cell style/number-format when a
# cell's type flips between text ("0" / "***.*" placeholders) and numeric:
#   C23 -> style for text "0"      E23 -> style for text "***.*"
#   J28 -> style for numeric count L14 -> style for numeric percent

# Row 15 (Rape)
$ws.Range("C23").Copy($ws.Range("C15"))
$ws.Range("J28").Copy($ws.Range("D15"))
$ws.Range("D15").Value = 1
$ws.Range("L14").Copy($ws.Range("E15"))
$ws.Range("E15").Value = -100
$ws.Range("J28").Copy($ws.Range("G15"))
$ws.Range("G15").Value = 1
$ws.Range("L14").Copy($ws.Range("H15"))
$ws.Range("H15").Value = 200
$ws.Range("I15").Value = 5
$ws.Range("J15").Value = 6
$ws.Range("K15").Value = -16.666666666666
$ws.Range("L15").Value = -16.666666666666
$ws.Range("N15").Value = 150

# Row 16 (Robbery)
$ws.Range("C16").Value = 5
$ws.Range("D16").Value = 2
$ws.Range("E16").Value = 150
$ws.Range("F16").Value = 10
$ws.Range("G16").Value = 16
$ws.Range("H16").Value = -37.5
$ws.Range("I16").Value = 37
$ws.Range("J16").Value = 42
$ws.Range("K16").Value = -11.904761904761
$ws.Range("L16").Value = 48
$ws.Range("M16").Value = 85
$ws.Range("N16").Value = -86.296296296296

# Row 17 (Fel. Assault)
$ws.Range("C17").Value = 3
$ws.Range("D17").Value = 4
$ws.Range("E17").Value = -25
$ws.Range("G17").Value = 10
$ws.Range("H17").Value = -10
$ws.Range("I17").Value = 34
$ws.Range("J17").Value = 35
$ws.Range("K17").Value = -2.857142857142
$ws.Range("L17").Value = -15
$ws.Range("M17").Value = 142.857142857143
$ws.Range("N17").Value = -38.181818181818

# Row 18 (Burglary)
$ws.Range("D18").Value = 4
$ws.Range("E18").Value = 0
$ws.Range("F18").Value = 17
$ws.Range("G18").Value = 23
$ws.Range("H18").Value = -26.086956521739
$ws.Range("I18").Value = 49
$ws.Range("J18").Value = 85
$ws.Range("K18").Value = -42.352941176470
$ws.Range("L18").Value = 44.117647058823
$ws.Range("M18").Value = -15.517241379310
$ws.Range("N18").Value = -80.933852140077

# Row 19 (Gr. Larceny)
$ws.Range("C19").Value = 31
$ws.Range("D19").Value = 21
$ws.Range("E19").Value = 47.619047619047
$ws.Range("F19").Value = 101
$ws.Range("G19").Value = 82
$ws.Range("H19").Value = 23.170731707317
$ws.Range("I19").Value = 327
$ws.Range("J19").Value = 324
$ws.Range("K19").Value = 0.925925925925
$ws.Range("L19").Value = 75.806451612903
$ws.Range("M19").Value = 3.481012658227
$ws.Range("N19").Value = -66.461538461538

# Row 20 (G.L.A.)
$ws.Range("C23").Copy($ws.Range("D20"))
$ws.Range("E23").Copy($ws.Range("E20"))
$ws.Range("F20").Value = 11
$ws.Range("H20").Value = 175
$ws.Range("I20").Value = 18
$ws.Range("K20").Value = 100
$ws.Range("L20").Value = 63.636363636363
$ws.Range("M20").Value = 350
$ws.Range("N20").Value = -92.653061224489

# Row 21 (TOTAL)
$ws.Range("C21").Value = 47
$ws.Range("D21").Value = 32
$ws.Range("E21").Value = 46.875
$ws.Range("F21").Value = 151
$ws.Range("G21").Value = 136
$ws.Range("H21").Value = 11.029411764705
$ws.Range("I21").Value = 470
$ws.Range("J21").Value = 501
$ws.Range("K21").Value = -6.187624750499
$ws.Range("L21").Value = 55.115511551155
$ws.Range("M21").Value = 14.077669902912
$ws.Range("N21").Value = -74.047487575924

# Row 22 (Transit)
$ws.Range("D22").Value = 2
$ws.Range("E22").Value = 0
$ws.Range("F22").Value = 4
$ws.Range("G22").Value = 12
$ws.Range("H22").Value = -66.666666666666
$ws.Range("I22").Value = 21
$ws.Range("J22").Value = 31
$ws.Range("K22").Value = -32.258064516129
$ws.Range("L22").Value = 0
$ws.Range("M22").Value = -8.695652173913

# Row 24 (Petit Larceny)
$ws.Range("C24").Value = 72
$ws.Range("D24").Value = 80
$ws.Range("E24").Value = -10
$ws.Range("F24").Value = 271
$ws.Range("G24").Value = 280
$ws.Range("H24").Value = -3.214285714285
$ws.Range("I24").Value = 1084
$ws.Range("J24").Value = 1061
$ws.Range("K24").Value = 2.167766258246
$ws.Range("L24").Value = 102.616822429907
$ws.Range("M24").Value = 136.681222707424

# Row 25 (Misd. Assault)
$ws.Range("C25").Value = 12
$ws.Range("D25").Value = 3
$ws.Range("E25").Value = 300
$ws.Range("F25").Value = 25
$ws.Range("G25").Value = 26
$ws.Range("H25").Value = -3.846153846153
$ws.Range("I25").Value = 90
$ws.Range("J25").Value = 87
$ws.Range("K25").Value = 3.448275862068
$ws.Range("L25").Value = 11.111111111111
$ws.Range("M25").Value = 32.352941176470

# Row 26 (UCR Rape*)
$ws.Range("C26").Value = 1
$ws.Range("J28").Copy($ws.Range("D26"))
$ws.Range("D26").Value = 1
$ws.Range("L14").Copy($ws.Range("E26"))
$ws.Range("E26").Value = 0
$ws.Range("F26").Value = 4
$ws.Range("J28").Copy($ws.Range("G26"))
$ws.Range("G26").Value = 1
$ws.Range("L14").Copy($ws.Range("H26"))
$ws.Range("H26").Value = 300
$ws.Range("I26").Value = 6
$ws.Range("J26").Value = 7
$ws.Range("K26").Value = -14.285714285714
$ws.Range("L26").Value = -14.285714285714

# Row 27 (Other Sex Crimes)
$ws.Range("C27").Value = 1
$ws.Range("D27").Value = 4
$ws.Range("E27").Value = -75
$ws.Range("F27").Value = 6
$ws.Range("G27").Value = 9
$ws.Range("H27").Value = -33.333333333333
$ws.Range("I27").Value = 23
$ws.Range("J27").Value = 24
$ws.Range("K27").Value = -4.166666666666
$ws.Range("L27").Value = 27.777777777777

# Row 30 (Hate Crimes)
$ws.Range("C23").Copy($ws.Range("G30"))
$ws.Range("E23").Copy($ws.Range("H30"))
